$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("correct")
$ws.Range("C2").Value = 1612
$ws.Range("D2").Value = 0.442
$ws.Range("C3").Value = 926
$ws.Range("D3").Value = 0.52
$ws.Range("C4").Value = 256
$ws.Range("D4").Value = 0.411
$ws.Range("C5").Value = 499
$ws.Range("D5").Value = 0.428
$ws.Range("C6").Value = 475
$ws.Range("D6").Value = 0.419
$ws.Range("C7").Value = 420
$ws.Range("D7").Value = 0.501
$ws.Range("C8").Value = 404
$ws.Range("D8").Value = 0.478
$ws.Range("C9").Value = 431
$ws.Range("D9").Value = 0.51
$ws.Range("C10").Value = 356
$ws.Range("D10").Value = 0.516
$ws.Range("C11").Value = 318
$ws.Range("D11").Value = 0.513
$ws.Range("C12").Value = 270
$ws.Range("D12").Value = 0.491
$ws.Range("C13").Value = 206
$ws.Range("D13").Value = 0.485
$ws.Range("C14").Value = 266
$ws.Range("D14").Value = 0.508
$ws.Range("C15").Value = 254
$ws.Range("D15").Value = 0.484
$ws.Range("C16").Value = 179
$ws.Range("D16").Value = 0.536
$ws.Range("C17").Value = 156
$ws.Range("D17").Value = 0.496

$ws = $wb.Worksheets.Item("distractor")
$ws.Range("C2").Value = 122
$ws.Range("D2").Value = -0.227
$ws.Range("C3").Value = 130
$ws.Range("D3").Value = -0.252
$ws.Range("C4").Value = 121
$ws.Range("D4").Value = -0.233
$ws.Range("C5").Value = 337
$ws.Range("D5").Value = -0.298
$ws.Range("C6").Value = 299
$ws.Range("D6").Value = -0.186
$ws.Range("C7").Value = 300
$ws.Range("D7").Value = -0.209
$ws.Range("C8").Value = 357
$ws.Range("D8").Value = -0.161
$ws.Range("C9").Value = 390
$ws.Range("D9").Value = -0.115
$ws.Range("C10").Value = 376
$ws.Range("D10").Value = -0.085
$ws.Range("D11").Value = -0.237
$ws.Range("C12").Value = 50
$ws.Range("D12").Value = -0.238
$ws.Range("D13").Value = -0.203
$ws.Range("C14").Value = 55
$ws.Range("D14").Value = -0.192
$ws.Range("C15").Value = 60
$ws.Range("D15").Value = -0.24
$ws.Range("C16").Value = 65
$ws.Range("D16").Value = -0.216
$ws.Range("C17").Value = 84
$ws.Range("D17").Value = -0.302
$ws.Range("C18").Value = 63
$ws.Range("D18").Value = -0.225
$ws.Range("C19").Value = 80
$ws.Range("D19").Value = -0.214
$ws.Range("C20").Value = 90
$ws.Range("D20").Value = -0.219
$ws.Range("C21").Value = 85
$ws.Range("D21").Value = -0.245
$ws.Range("C22").Value = 69
$ws.Range("D22").Value = -0.238
$ws.Range("C23").Value = 81
$ws.Range("D23").Value = -0.266
$ws.Range("C24").Value = 79
$ws.Range("D24").Value = -0.238
$ws.Range("C25").Value = 85
$ws.Range("D25").Value = -0.249
$ws.Range("C26").Value = 90
$ws.Range("D26").Value = -0.236
$ws.Range("C27").Value = 113
$ws.Range("D27").Value = -0.254
$ws.Range("C28").Value = 84
$ws.Range("D28").Value = -0.232
$ws.Range("C29").Value = 109
$ws.Range("D29").Value = -0.23
$ws.Range("C30").Value = 95
$ws.Range("D30").Value = -0.199
$ws.Range("C31").Value = 101
$ws.Range("D31").Value = -0.265
$ws.Range("C32").Value = 116
$ws.Range("D32").Value = -0.227
$ws.Range("C33").Value = 80
$ws.Range("D33").Value = -0.166
$ws.Range("D34").Value = -0.25
$ws.Range("C35").Value = 98
$ws.Range("D35").Value = -0.186
$ws.Range("C36").Value = 93
$ws.Range("D36").Value = -0.193
$ws.Range("C37").Value = 88
$ws.Range("D37").Value = -0.231
$ws.Range("C38").Value = 122
$ws.Range("D38").Value = -0.166
$ws.Range("D39").Value = -0.196
$ws.Range("C40").Value = 151
$ws.Range("D40").Value = -0.256
$ws.Range("C41").Value = 139
$ws.Range("D41").Value = -0.174
$ws.Range("C42").Value = 138
$ws.Range("D42").Value = -0.171
$ws.Range("C43").Value = 147
$ws.Range("D43").Value = -0.232
$ws.Range("C44").Value = 159
$ws.Range("D44").Value = -0.124
$ws.Range("C45").Value = 186
$ws.Range("D45").Value = -0.186
$ws.Range("D46").Value = -0.241
$ws.Range("C47").Value = 167
$ws.Range("D47").Value = -0.163
$ws.Range("C48").Value = 179
$ws.Range("D48").Value = -0.164
$ws.Range("C49").Value = 174
$ws.Range("D49").Value = -0.151

$ws = $wb.Worksheets.Item("descriptives")
$ws.Range("B2").Value = 0.484
$ws.Range("C2").Value = -0.212
$ws.Range("B3").Value = 0.038
$ws.Range("C3").Value = 0.044
$ws.Range("B4").Value = 0.494
$ws.Range("C4").Value = -0.226
$ws.Range("B5").Value = 0.411
$ws.Range("C5").Value = -0.302
$ws.Range("B6").Value = 0.536
$ws.Range("C6").Value = -0.085
